$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 5 and 4 (bottom-up to keep row indices stable)
$ws.Rows.Item(5).Delete()
$ws.Rows.Item(4).Delete()

# --- Row 2 updates ---
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "1"
$ws.Range("D2").Value = -0.0415
$ws.Range("E2").Value = -0.115
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 9.99
$ws.Range("L2").Value = 0.15
$ws.Range("M2").Value = 9.300000000000001
$ws.Range("N2").Value = 0.2319201995012469
$ws.Range("O2").Value = 0.9309309309309309
$ws.Range("P2").Value = 9.300000000000001
$ws.Range("Q2").Value = 0.2319201995012469
$ws.Range("R2").Value = 0.9309309309309309
$ws.Range("U2").Value = 280.4
$ws.Range("V2").Value = 6.992518703241895
$ws.Range("W2").Value = 0.09541547277936963
$ws.Range("X2").Value = 0.09010471541889752
$ws.Range("Y2").Value = 0.00531075736047211
$ws.Range("Z2").Value = -0.4078383343539497
$ws.Range("AA2").Value = -0
$ws.Range("AB2").Value = 0.04605330149756055
$ws.Range("AC2").Value = -0.04605330149756055
$ws.Range("AD2").Value = 105
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 105
$ws.Range("AG2").Value = -175.4
$ws.Range("AH2").Value = 0.7236388697450035
$ws.Range("AI2").Value = 0.523168908819133
$ws.Range("AJ2").Value = 1.296378418329638
$ws.Range("AK2").Value = 2.200752823086575
$ws.Range("AN2").ClearContents()
$ws.Range("AP2").ClearContents()

# --- Row 3 updates ---
$ws.Range("B3").Value = "Standard Chartered Bank Botswana Limited (BSM:STANCHART)"
$ws.Range("D3").Value = -0.0415
$ws.Range("E3").Value = -0.115
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 9.99
$ws.Range("L3").Value = 0.15
$ws.Range("M3").Value = 9.300000000000001
$ws.Range("N3").Value = 0.2319201995012469
$ws.Range("O3").Value = 0.9309309309309309
$ws.Range("P3").Value = 9.300000000000001
$ws.Range("Q3").Value = 0.2319201995012469
$ws.Range("R3").Value = 0.9309309309309309
$ws.Range("U3").Value = 280.4
$ws.Range("V3").Value = 6.992518703241895
$ws.Range("W3").Value = 0.09541547277936963
$ws.Range("X3").Value = 0.09010471541889752
$ws.Range("Y3").Value = 0.00531075736047211
$ws.Range("Z3").Value = -0.4078383343539497
$ws.Range("AA3").Value = -0
$ws.Range("AB3").Value = 0.04605330149756055
$ws.Range("AC3").Value = -0.04605330149756055
$ws.Range("AD3").Value = 105
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 105
$ws.Range("AG3").Value = -175.4
$ws.Range("AH3").Value = 0.7236388697450035
$ws.Range("AI3").Value = 0.523168908819133
$ws.Range("AJ3").Value = 1.296378418329638
$ws.Range("AK3").Value = 2.200752823086575
$ws.Range("AN3").ClearContents()
$ws.Range("AP3").ClearContents()
